$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# Clear the word-pair text that was in B1:D15 (memory-key words), keeping
# the existing cell formatting/styles intact.
$ws.Range("B1:D15").ClearContents()

# Reset the scrolled viewport back to the top-left (drops topLeftCell="A4").
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1

# Reproduce the final selection state: cursor on I13, then Select All.
$ws.Cells.Item(13, 9).Select()
$ws.Cells.Select()
